# [Fonds de solidarite] Add 2020-12-24 data
# Updates nombre_aides (col C) and montant_total (col D) for the rows whose
# underlying counts/amounts changed with the new data refresh.
# Values are stored as text in this workbook (matching the original inlineStr
# cells), so a leading apostrophe is used to stop Excel's COM layer from
# auto-converting the numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row = 3; C = "1412"; D = "6930416.29"}
    @{Row = 4; C = "624"; D = "5729846.20"}
    @{Row = 6; C = "72"; D = "1191921.05"}
    @{Row = 9; C = "89"; D = "327730.09"}
    @{Row = 10; C = "480"; D = "2604112.26"}
    @{Row = 12; C = "62"; D = "774841.00"}
    @{Row = 13; C = "19"; D = "398076.00"}
    @{Row = 14; C = "30"; D = "83163.21"}
    @{Row = 15; C = "137"; D = "763392.38"}
    @{Row = 16; C = "526"; D = "2650081.59"}
    @{Row = 23; C = "407"; D = "1958207.73"}
    @{Row = 24; C = "165"; D = "1308657.40"}
    @{Row = 27; C = "8"; D = "227000.00"}
    @{Row = 35; C = "159"; D = "653861.72"}
    @{Row = 36; C = "781"; D = "3562981.52"}
    @{Row = 37; C = "376"; D = "3139135.95"}
    @{Row = 38; C = "138"; D = "1689049.71"}
    @{Row = 39; C = "56"; D = "1027468.00"}
    @{Row = 49; C = "173"; D = "821969.56"}
    @{Row = 50; C = "93"; D = "639868.00"}
    @{Row = 51; C = "41"; D = "340681.00"}
    @{Row = 53; C = "3"; D = "65000.00"}
    @{Row = 56; C = "1013"; D = "5645922.71"}
    @{Row = 57; C = "488"; D = "4369471.82"}
    @{Row = 58; C = "177"; D = "1707927.84"}
    @{Row = 60; C = "13"; D = "261032.25"}
    @{Row = 64; C = "3112"; D = "18552283.83"}
    @{Row = 65; C = "1088"; D = "8079411.26"}
    @{Row = 83; C = "119"; D = "397139.87"}
    @{Row = 84; C = "499"; D = "2331055.30"}
    @{Row = 85; C = "203"; D = "1643538.92"}
    @{Row = 86; C = "80"; D = "1055060.51"}
    @{Row = 87; C = "27"; D = "433568.04"}
    @{Row = 90; C = "287"; D = "1293756.14"}
    @{Row = 91; C = "1094"; D = "5670038.26"}
    @{Row = 92; C = "461"; D = "4036769.31"}
    @{Row = 93; C = "187"; D = "2169030.63"}
    @{Row = 94; C = "67"; D = "1264435.17"}
    @{Row = 95; C = "11"; D = "425000.00"}
    @{Row = 96; C = "42"; D = "144991.00"}
    @{Row = 97; C = "385"; D = "1556656.97"}
    @{Row = 98; C = "1343"; D = "6344691.79"}
    @{Row = 99; C = "516"; D = "3564002.92"}
    @{Row = 100; C = "179"; D = "1900457.65"}
    @{Row = 101; C = "56"; D = "867508.00"}
    @{Row = 102; C = "7"; D = "200000.00"}
    @{Row = 110; C = "441"; D = "1584017.20"}
    @{Row = 111; C = "1766"; D = "7496015.91"}
    @{Row = 112; C = "726"; D = "4798137.64"}
    @{Row = 113; C = "244"; D = "2933638.27"}
    @{Row = 114; C = "81"; D = "1126148.00"}
    @{Row = 115; C = "9"; D = "228173.00"}
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = "'" + $u.C
    $ws.Range("D" + $u.Row).Value = "'" + $u.D
}
